$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1758.1
$ws.Range("I19").Value = 1984.625
$ws.Range("J19").Value = 852
$ws.Range("K19").Value = 1984.625
$ws.Range("L19").Value = 852
$ws.Range("M19").Value = -1809.625
$ws.Range("N19").Value = -1202
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H41").Value = 447.5
$ws.Range("I41").Value = 447.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 447.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -7.5
$ws.Range("N41").ClearContents()
$ws.Range("H64").Value = 3499.5
$ws.Range("I64").Value = 3499.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3499.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3251.5
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3499.5
$ws.Range("I67").Value = 3499.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3499.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2641.5
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 1587.6666
$ws.Range("I100").Value = 1587.6666
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1587.6666
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1046.6666
$ws.Range("N100").ClearContents()
$ws.Range("H116").Value = 3401
$ws.Range("I116").Value = 2835
$ws.Range("J116").Value = 4250
$ws.Range("K116").Value = 2835
$ws.Range("L116").Value = 4250
$ws.Range("M116").Value = 607
$ws.Range("N116").Value = -11134

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3399.5
$ws.Range("I61").Value = 3399.5
$ws.Range("K61").Value = 3399.5
$ws.Range("M61").Value = -3187.5
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
$ws.Range("H136").Value = 3399.5
$ws.Range("I136").Value = 3399.5
$ws.Range("K136").Value = 10198.5
$ws.Range("M136").Value = -7648.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 74998
$ws.Range("J13").Value = 74998
$ws.Range("L13").Value = 74998
$ws.Range("N13").Value = -75334
$ws.Range("H76").Value = 21964.666
$ws.Range("J76").Value = 21964.666
$ws.Range("L76").Value = 21964.666
$ws.Range("N76").Value = -22594.666
$ws.Range("H79").Value = 21964.666
$ws.Range("J79").Value = 21964.666
$ws.Range("L79").Value = 21964.666
$ws.Range("N79").Value = -24148.666
$ws.Range("H105").Value = 9749.25
$ws.Range("I105").Value = 8999
$ws.Range("K105").Value = 8999
$ws.Range("M105").Value = -7252

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H93").Value = 29866.166
$ws.Range("J93").Value = 34999.75
$ws.Range("L93").Value = 34999.75
$ws.Range("N93").Value = -38743.75
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H107").Value = 3799.3333
$ws.Range("I107").Value = 3799.3333
$ws.Range("K107").Value = 3799.3333
$ws.Range("M107").Value = -1879.3333
$ws.Range("H116").Value = 40742
$ws.Range("J116").Value = 40742
$ws.Range("L116").Value = 40742
$ws.Range("N116").Value = -49920
$ws.Range("H122").Value = 2413.3333
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 2370
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 7110
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -12010
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H137").Value = 89831.664
$ws.Range("J137").Value = 89748.5
$ws.Range("L137").Value = 89748.5
$ws.Range("N137").Value = -99948.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 598.3333
$ws.Range("I2").Value = 286.6
$ws.Range("J2").Value = 2157
$ws.Range("K2").Value = 1719.6
$ws.Range("L2").Value = 12942
$ws.Range("M2").Value = -1606.6
$ws.Range("N2").Value = -13168
$ws.Range("H23").Value = 206.88889
$ws.Range("J23").Value = 172.2
$ws.Range("L23").Value = 516.5999999999999
$ws.Range("N23").Value = -986.5999999999999
$ws.Range("H99").Value = 2219.6
$ws.Range("I99").Value = 1524.5
$ws.Range("K99").Value = 4573.5
$ws.Range("M99").Value = -2327.5
$ws.Range("H129").Value = 2823.6667
$ws.Range("I129").Value = 2940
$ws.Range("J129").Value = 2765.5
$ws.Range("K129").Value = 8820
$ws.Range("L129").Value = 8296.5
$ws.Range("M129").Value = -3820
$ws.Range("N129").Value = -18296.5
$ws.Range("H131").Value = 4113.7144
$ws.Range("I131").Value = 14999
$ws.Range("J131").Value = 2299.5
$ws.Range("K131").Value = 44997
$ws.Range("L131").Value = 6898.5
$ws.Range("M131").Value = -39957
$ws.Range("N131").Value = -16978.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 482.8
$ws.Range("I97").Value = 532.5
$ws.Range("J97").Value = 449.66666
$ws.Range("K97").Value = 532.5
$ws.Range("L97").Value = 449.66666
$ws.Range("M97").Value = -36.5
$ws.Range("N97").Value = -1441.66666
$ws.Range("H132").Value = 5988.1113
$ws.Range("I132").Value = 7069.8335
$ws.Range("J132").Value = 3824.6667
$ws.Range("K132").Value = 21209.5005
$ws.Range("L132").Value = 11474.0001
$ws.Range("M132").Value = -18679.5005
$ws.Range("N132").Value = -16534.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3374.5
$ws.Range("I93").Value = 3374.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3374.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2126.5
$ws.Range("N93").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 4536.3335
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 4536.3335
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4536.3335
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -5518.3335
$ws.Range("H107").Value = 1700.7142
$ws.Range("I107").Value = 1579.8
$ws.Range("K107").Value = 4739.4
$ws.Range("M107").Value = -2819.4
$ws.Range("H136").Value = 13953.077
$ws.Range("I136").Value = 14717.272
$ws.Range("J136").Value = 9750
$ws.Range("K136").Value = 44151.81600000001
$ws.Range("L136").Value = 29250
$ws.Range("M136").Value = -41601.81600000001
$ws.Range("N136").Value = -34350
